$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3364.375
$ws.Range("I2").Value = 786
$ws.Range("J2").Value = 7661.6665
$ws.Range("K2").Value = 786
$ws.Range("L2").Value = 7661.6665
$ws.Range("M2").Value = -673
$ws.Range("N2").Value = -7887.6665
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H53").Value = 292.72726
$ws.Range("I53").Value = 239.64706
$ws.Range("K53").Value = 239.64706
$ws.Range("M53").Value = 397.35294
$ws.Range("H86").Value = 6042.091
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 6246.3
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 6246.3
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -8492.299999999999
$ws.Range("H89").Value = 6042.091
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 6246.3
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 31231.5
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -42463.5
$ws.Range("H92").Value = 447
$ws.Range("I92").Value = 412.15384
$ws.Range("K92").Value = 412.15384
$ws.Range("M92").Value = 835.8461600000001
$ws.Range("H106").Value = 1994
$ws.Range("I106").Value = 1994
$ws.Range("K106").Value = 1994
$ws.Range("M106").Value = -1363
$ws.Range("H116").Value = 5401.5
$ws.Range("I116").Value = 2803
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 2803
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = 639
$ws.Range("N116").Value = -14884
$ws.Range("H132").Value = 1338.4166
$ws.Range("I132").Value = 1369.1818
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 4107.5454
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1577.5454
$ws.Range("N132").Value = -8060
$ws.Range("H138").Value = 3501.987
$ws.Range("J138").Value = 3419.0156
$ws.Range("L138").Value = 10257.0468
$ws.Range("N138").Value = -20537.0468

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1753.5264
$ws.Range("I2").Value = 811.36365
$ws.Range("J2").Value = 3049
$ws.Range("K2").Value = 811.36365
$ws.Range("L2").Value = 3049
$ws.Range("M2").Value = -698.36365
$ws.Range("N2").Value = -3275
$ws.Range("H61").Value = 5380.857
$ws.Range("I61").Value = 1890
$ws.Range("K61").Value = 1890
$ws.Range("M61").Value = -1678
$ws.Range("H74").Value = 2572.0557
$ws.Range("I74").Value = 2356.1428
$ws.Range("K74").Value = 2356.1428
$ws.Range("M74").Value = -1482.1428
$ws.Range("H77").Value = 2572.0557
$ws.Range("I77").Value = 2356.1428
$ws.Range("K77").Value = 11780.714
$ws.Range("M77").Value = -7412.714
$ws.Range("H116").Value = 1753.5264
$ws.Range("I116").Value = 811.36365
$ws.Range("J116").Value = 3049
$ws.Range("K116").Value = 811.36365
$ws.Range("L116").Value = 3049
$ws.Range("M116").Value = 1482.63635
$ws.Range("N116").Value = -7637
$ws.Range("H136").Value = 5380.857
$ws.Range("I136").Value = 1890
$ws.Range("K136").Value = 5670
$ws.Range("M136").Value = -3120

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1753.5264
$ws.Range("I3").Value = 811.36365
$ws.Range("J3").Value = 3049
$ws.Range("K3").Value = 811.36365
$ws.Range("L3").Value = 3049
$ws.Range("M3").Value = -697.36365
$ws.Range("N3").Value = -3277
$ws.Range("H107").Value = 2909.1177
$ws.Range("I107").Value = 2590.125
$ws.Range("K107").Value = 2590.125
$ws.Range("M107").Value = -670.125
$ws.Range("H134").Value = 799.5
$ws.Range("I134").Value = 799.5
$ws.Range("K134").Value = 2398.5
$ws.Range("M134").Value = 136.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3520.4614
$ws.Range("I31").Value = 3225
$ws.Range("K31").Value = 3225
$ws.Range("M31").Value = -2930
$ws.Range("H34").Value = 3520.4614
$ws.Range("I34").Value = 3225
$ws.Range("K34").Value = 3225
$ws.Range("M34").Value = -3023
$ws.Range("H134").Value = 2133.32
$ws.Range("I134").Value = 2063.9583
$ws.Range("J134").Value = 3798
$ws.Range("K134").Value = 6191.874899999999
$ws.Range("L134").Value = 11394
$ws.Range("M134").Value = -3656.874899999999
$ws.Range("N134").Value = -16464

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 497.875
$ws.Range("I122").Value = 515.5
$ws.Range("J122").Value = 492
$ws.Range("K122").Value = 4639.5
$ws.Range("L122").Value = 4428
$ws.Range("M122").Value = -2189.5
$ws.Range("N122").Value = -9328

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4723.7
$ws.Range("I70").Value = 4627.5713
$ws.Range("K70").Value = 4627.5713
$ws.Range("M70").Value = -4357.5713
$ws.Range("H73").Value = 4723.7
$ws.Range("I73").Value = 4627.5713
$ws.Range("K73").Value = 4627.5713
$ws.Range("M73").Value = -3691.5713
$ws.Range("H122").Value = 4256.421
$ws.Range("I122").Value = 3490.6428
$ws.Range("J122").Value = 6400.6
$ws.Range("K122").Value = 10471.9284
$ws.Range("L122").Value = 19201.8
$ws.Range("M122").Value = -8021.928400000001
$ws.Range("N122").Value = -24101.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6933.607
$ws.Range("I7").Value = 3750.1428
$ws.Range("J7").Value = 7994.7617
$ws.Range("K7").Value = 3750.1428
$ws.Range("L7").Value = 7994.7617
$ws.Range("M7").Value = -3638.1428
$ws.Range("N7").Value = -8218.761699999999
$ws.Range("H40").Value = 3752.2666
$ws.Range("I40").Value = 2869.2856
$ws.Range("J40").Value = 4524.875
$ws.Range("K40").Value = 2869.2856
$ws.Range("L40").Value = 4524.875
$ws.Range("M40").Value = -2733.2856
$ws.Range("N40").Value = -4796.875
$ws.Range("H122").Value = 5264.1177
$ws.Range("I122").Value = 4560.6
$ws.Range("K122").Value = 13681.8
$ws.Range("M122").Value = -11231.8
$ws.Range("H126").Value = 6933.607
$ws.Range("I126").Value = 3750.1428
$ws.Range("J126").Value = 7994.7617
$ws.Range("K126").Value = 11250.4284
$ws.Range("L126").Value = 23984.2851
$ws.Range("M126").Value = -8780.428400000001
$ws.Range("N126").Value = -28924.2851

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3329.3
$ws.Range("I132").Value = 3049.8333
$ws.Range("J132").Value = 3748.5
$ws.Range("K132").Value = 9149.499899999999
$ws.Range("L132").Value = 11245.5
$ws.Range("M132").Value = -6619.499899999999
$ws.Range("N132").Value = -16305.5
$ws.Range("H136").Value = 4846.8335
